$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows are rotated: the pair of rows that used to hold the
# most recent date (rows 6-7) now becomes the first pair (rows 2-3); the
# pair that used to be first (rows 2-3) moves down to rows 4-5; and the
# middle pair (rows 4-5) moves down to rows 6-7.

# Row 2 (was the old row 6 data)
$ws.Range("D2").Value = 44574
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6500
$ws.Range("S2").Value = 3250

# Row 3 (was the old row 7 data)
$ws.Range("D3").Value = 44574
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 5000
$ws.Range("O3").Value = 5000
$ws.Range("P3").Value = 5000
$ws.Range("S3").Value = 2500

# Row 4 (was the old row 2 data)
$ws.Range("D4").Value = 44223
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 3500
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3750
$ws.Range("S4").Value = 1875

# Row 5 (was the old row 3 data)
$ws.Range("D5").Value = 44223
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 3000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("S5").Value = 1500

# Row 6 (was the old row 4 data)
$ws.Range("D6").Value = 44559

# Row 7 (was the old row 5 data)
$ws.Range("D7").Value = 44559
